# Applies the cryptos.xlsx update (2023-02-15 09:16:30 UTC GitHub Actions run):
# refreshed coin prices / 1h volume %, rotated a few coin listings, and bumped the "Hora" column from 8 to 9.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. Every value must stay plain text
# (matching the original inlineStr cells), so NumberFormat is forced to "@"
# for each cell before the value is assigned - this stops Excel from turning
# things like "297.90" or "1.99%" into numeric/percentage cells.
$updates = [ordered]@{
    'D2' = '297.90'
    'E2' = '1.99%'
    'G2' = '9'
    'D3' = '42.07'
    'E3' = '4.48%'
    'G3' = '9'
    'D4' = '5.010'
    'E4' = '-0.05%'
    'G4' = '9'
    'D5' = '0.07524'
    'E5' = '3.28%'
    'G5' = '9'
    'B6' = 'FTXToken'
    'C6' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D6' = '1.595'
    'E6' = '3.60%'
    'G6' = '9'
    'B7' = 'MXToken'
    'C7' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D7' = '0.9172'
    'E7' = '-0.60%'
    'G7' = '9'
    'B8' = 'BTSEToken'
    'C8' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D8' = '2.401'
    'E8' = '2.14%'
    'G8' = '9'
    'B9' = 'LiechtensteinCryptoassetsExchange'
    'C9' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D9' = '0.1185'
    'E9' = '2.42%'
    'G9' = '9'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D10' = '0.1830'
    'E10' = '4.61%'
    'G10' = '9'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D11' = '0.08940'
    'E11' = '2.63%'
    'G11' = '9'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D12' = '0.04129'
    'E12' = '-5.26%'
    'G12' = '9'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D13' = '0.1050'
    'E13' = '-0.23%'
    'G13' = '9'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D14' = '0.001277'
    'E14' = '0.17%'
    'G14' = '9'
    'B15' = 'TigerCash'
    'C15' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D15' = '0.006010'
    'E15' = '-0.47%'
    'G15' = '9'
    'B16' = 'LEO'
    'C16' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D16' = '3.343'
    'E16' = '0.14%'
    'G16' = '9'
    'B17' = 'GateToken'
    'C17' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D17' = '4.368'
    'E17' = '1.77%'
    'G17' = '9'
    'D18' = '0.3331'
    'E18' = '1.39%'
    'G18' = '9'
    'D19' = '8.368'
    'E19' = '6.32%'
    'G19' = '9'
    'E20' = '-2.86%'
    'G20' = '9'
    'E21' = '11.94%'
    'G21' = '9'
    'D22' = '0.04093'
    'E22' = '4.51%'
    'G22' = '9'
    'E23' = '0.30%'
    'G23' = '9'
    'D24' = '0.003891'
    'E24' = '6.63%'
    'G24' = '9'
    'E25' = '8.36%'
    'G25' = '9'
    'G26' = '9'
    'G27' = '9'
    'G28' = '9'
    'G29' = '9'
    'G30' = '9'
    'G31' = '9'
    'G32' = '9'
    'G33' = '9'
    'G34' = '9'
    'G35' = '9'
    'G36' = '9'
    'G37' = '9'
    'D38' = '0.02388'
    'E38' = '4.04%'
    'G38' = '9'
    'D39' = '0.05230'
    'E39' = '3.68%'
    'G39' = '9'
    'D40' = '0.007042'
    'E40' = '27.50%'
    'G40' = '9'
    'D41' = '0.007786'
    'E41' = '-0.75%'
    'G41' = '9'
    'D42' = '0.1327'
    'E42' = '3.21%'
    'G42' = '9'
    'D43' = '0.007403'
    'E43' = '0.33%'
    'G43' = '9'
    'D44' = '0.007128'
    'E44' = '-1.99%'
    'G44' = '9'
    'D45' = '0.3222'
    'E45' = '1.20%'
    'G45' = '9'
    'D46' = '0.00006565'
    'E46' = '6.09%'
    'G46' = '9'
    'E47' = '-0.08%'
    'G47' = '9'
    'D48' = '0.04539'
    'E48' = '-12.45%'
    'G48' = '9'
    'D49' = '0.004207'
    'E49' = '0.12%'
    'G49' = '9'
    'E50' = '-0.08%'
    'G50' = '9'
    'E51' = '-0.08%'
    'G51' = '9'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
